$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns F, G, H in row 155 (second table header)
$ws.Range("F155").Value = "Sum Lines added"
$ws.Range("G155").Value = " Sum Lines removed"
$ws.Range("H155").Value = " Sum No. Commits"

# Data rows 156-226: F = Lines added, G = Lines removed, H = No. Commits
$arr = New-Object "object[,]" 71,3
$arr[0,0] = 180; $arr[0,1] = 0; $arr[0,2] = 2
$arr[1,0] = 16; $arr[1,1] = 0; $arr[1,2] = 2
$arr[2,0] = 14; $arr[2,1] = 0; $arr[2,2] = 1
$arr[3,0] = 19; $arr[3,1] = 0; $arr[3,2] = 2
$arr[4,0] = 22; $arr[4,1] = 0; $arr[4,2] = 2
$arr[5,0] = 791; $arr[5,1] = 0; $arr[5,2] = 41
$arr[6,0] = 13; $arr[6,1] = 0; $arr[6,2] = 4
$arr[7,0] = 332; $arr[7,1] = 0; $arr[7,2] = 1
$arr[8,0] = 91; $arr[8,1] = 0; $arr[8,2] = 6
$arr[9,0] = 43; $arr[9,1] = 0; $arr[9,2] = 1
$arr[10,0] = 799; $arr[10,1] = 0; $arr[10,2] = 41
$arr[11,0] = 333; $arr[11,1] = 0; $arr[11,2] = 1
$arr[12,0] = 10; $arr[12,1] = 0; $arr[12,2] = 1
$arr[13,0] = 299; $arr[13,1] = 0; $arr[13,2] = 85
$arr[14,0] = 420; $arr[14,1] = 0; $arr[14,2] = 24
$arr[15,0] = 36; $arr[15,1] = 0; $arr[15,2] = 6
$arr[16,0] = 11; $arr[16,1] = 0; $arr[16,2] = 3
$arr[17,0] = 110; $arr[17,1] = 0; $arr[17,2] = 3
$arr[18,0] = 189; $arr[18,1] = 0; $arr[18,2] = 29
$arr[19,0] = 32; $arr[19,1] = 0; $arr[19,2] = 2
$arr[20,0] = 83; $arr[20,1] = 0; $arr[20,2] = 3
$arr[21,0] = 101; $arr[21,1] = 0; $arr[21,2] = 8
$arr[22,0] = 15; $arr[22,1] = 0; $arr[22,2] = 10
$arr[23,0] = 484; $arr[23,1] = 0; $arr[23,2] = 11
$arr[24,0] = 237; $arr[24,1] = 0; $arr[24,2] = 1
$arr[25,0] = 376; $arr[25,1] = 0; $arr[25,2] = 1
$arr[26,0] = 68; $arr[26,1] = 0; $arr[26,2] = 1
$arr[27,0] = 230; $arr[27,1] = 0; $arr[27,2] = 2
$arr[28,0] = 124; $arr[28,1] = 0; $arr[28,2] = 20
$arr[29,0] = 66; $arr[29,1] = 0; $arr[29,2] = 7
$arr[30,0] = 95; $arr[30,1] = 0; $arr[30,2] = 12
$arr[31,0] = 16; $arr[31,1] = 0; $arr[31,2] = 1
$arr[32,0] = 153; $arr[32,1] = 0; $arr[32,2] = 7
$arr[33,0] = 211; $arr[33,1] = 0; $arr[33,2] = 4
$arr[34,0] = 6; $arr[34,1] = 0; $arr[34,2] = 1
$arr[35,0] = 121; $arr[35,1] = 0; $arr[35,2] = 7
$arr[36,0] = 198; $arr[36,1] = 0; $arr[36,2] = 3
$arr[37,0] = 37; $arr[37,1] = 0; $arr[37,2] = 3
$arr[38,0] = 82; $arr[38,1] = 0; $arr[38,2] = 11
$arr[39,0] = 96; $arr[39,1] = 0; $arr[39,2] = 9
$arr[40,0] = 55; $arr[40,1] = 0; $arr[40,2] = 3
$arr[41,0] = 16; $arr[41,1] = 0; $arr[41,2] = 1
$arr[42,0] = 5; $arr[42,1] = 0; $arr[42,2] = 2
$arr[43,0] = 1812; $arr[43,1] = 0; $arr[43,2] = 31
$arr[44,0] = 19; $arr[44,1] = 0; $arr[44,2] = 2
$arr[45,0] = 106; $arr[45,1] = 0; $arr[45,2] = 7
$arr[46,0] = 22; $arr[46,1] = 0; $arr[46,2] = 2
$arr[47,0] = 4; $arr[47,1] = 0; $arr[47,2] = 2
$arr[48,0] = 70; $arr[48,1] = 0; $arr[48,2] = 24
$arr[49,0] = 789; $arr[49,1] = 0; $arr[49,2] = 4
$arr[50,0] = 66; $arr[50,1] = 0; $arr[50,2] = 17
$arr[51,0] = 145; $arr[51,1] = 0; $arr[51,2] = 3
$arr[52,0] = 105; $arr[52,1] = 0; $arr[52,2] = 8
$arr[53,0] = 32; $arr[53,1] = 0; $arr[53,2] = 4
$arr[54,0] = 59; $arr[54,1] = 0; $arr[54,2] = 6
$arr[55,0] = 290; $arr[55,1] = 0; $arr[55,2] = 3
$arr[56,0] = 40; $arr[56,1] = 0; $arr[56,2] = 1
$arr[57,0] = 20; $arr[57,1] = 0; $arr[57,2] = 1
$arr[58,0] = 17; $arr[58,1] = 0; $arr[58,2] = 1
$arr[59,0] = 38; $arr[59,1] = 0; $arr[59,2] = 1
$arr[60,0] = 51; $arr[60,1] = 0; $arr[60,2] = 11
$arr[61,0] = 538; $arr[61,1] = 0; $arr[61,2] = 17
$arr[62,0] = 10; $arr[62,1] = 0; $arr[62,2] = 1
$arr[63,0] = 14; $arr[63,1] = 0; $arr[63,2] = 7
$arr[64,0] = 107; $arr[64,1] = 0; $arr[64,2] = 12
$arr[65,0] = 32; $arr[65,1] = 0; $arr[65,2] = 2
$arr[66,0] = 81; $arr[66,1] = 0; $arr[66,2] = 4
$arr[67,0] = 313; $arr[67,1] = 0; $arr[67,2] = 3
$arr[68,0] = 501; $arr[68,1] = 0; $arr[68,2] = 36
$arr[69,0] = 95; $arr[69,1] = 0; $arr[69,2] = 5
$arr[70,0] = 16; $arr[70,1] = 0; $arr[70,2] = 10
$ws.Range("F156:H226").Value = $arr

# Update the active selection to match the edited workbook state
[void]$ws.Range("F156").Select()
